# Important tweak delta trap figure
#
# The two right-most "Delta" labels in the trap figure are renamed to
# "Modified" and widened (re-centered on their original midpoint) to fit
# the new, longer text.
#
# Shape.Left/.Top/.Width/.Height are expressed in points (1 pt = 12700 EMU).
# The point values below are chosen so that they reproduce the exact EMU
# geometry from the source slide (PowerPoint's Shape position/size
# properties are single-precision floats, so the literals are picked to
# land on the correct EMU value once converted/truncated).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 3rd "Delta" textbox (TextBox 38 / id 39)
$shape1 = $s.Shapes.Item(13)
$shape1.Left = 448.94193
$shape1.Top = 364.022245
$shape1.Width = 59.222875
$shape1.Height = 21.810985
$shape1.TextFrame.TextRange.Text = "Modified"

# 4th "Delta" textbox (TextBox 39 / id 40)
$shape2 = $s.Shapes.Item(14)
$shape2.Left = 500.429814
$shape2.Top = 364.022245
$shape2.Width = 59.222875
$shape2.Height = 21.810985
$shape2.TextFrame.TextRange.Text = "Modified"
